$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Datos actualizados" timestamp (A1)
$ws.Range("A1").Value = "Datos actualizados a 16 de Septiembre de 2020 a las 12:01"

# Swap Santa Lucia / Timor Oriental order (rows 204-205)
$ws.Range("A204").Value = "Timor Oriental"
$ws.Range("A205").Value = "Santa Lucia"

# Swap Montserrat / Islas Malvinas order (rows 214-215), including their stats
$ws.Range("A214").Value = "Islas Malvinas"
$ws.Range("A215").Value = "Montserrat"

# Update COVID-19 statistics for updated countries
# Row 5
$ws.Range("B5").Value = 5025043
$ws.Range("C5").Value = 7009
$ws.Range("E5").Value = 1000559
$ws.Range("G5").Value = 33
$ws.Range("H5").Value = 82124

# Row 7
$ws.Range("B7").Value = 1079519
$ws.Range("C7").Value = 5670
$ws.Range("D7").Value = 890114
$ws.Range("E7").Value = 170488
$ws.Range("G7").Value = 132
$ws.Range("H7").Value = 18917

# Row 18
$ws.Range("B18").Value = 342671
$ws.Range("C18").Value = 1615
$ws.Range("D18").Value = 247969
$ws.Range("E18").Value = 89879
$ws.Range("G18").Value = 21
$ws.Range("H18").Value = 4823

# Row 26
$ws.Range("B26").Value = 228993
$ws.Range("C26").Value = 3963
$ws.Range("D26").Value = 164101
$ws.Range("E26").Value = 55792
$ws.Range("G26").Value = 135
$ws.Range("H26").Value = 9100

# Row 40
$ws.Range("B40").Value = 91196
$ws.Range("C40").Value = 536
$ws.Range("D40").Value = 84363
$ws.Range("E40").Value = 6028
$ws.Range("G40").Value = 8
$ws.Range("H40").Value = 805

# Row 54
$ws.Range("E54").Value = 6597
$ws.Range("G54").Value = 2
$ws.Range("H54").Value = 215

# Row 70
$ws.Range("B70").Value = 35073
$ws.Range("C70").Value = 768
$ws.Range("D70").Value = 27655
$ws.Range("E70").Value = 6660
$ws.Range("G70").Value = 1
$ws.Range("H70").Value = 758

# Row 95
$ws.Range("B95").Value = 10414
$ws.Range("C95").Value = 13
$ws.Range("D95").Value = 9827
$ws.Range("E95").Value = 320

# Row 97
$ws.Range("B97").Value = 10031
$ws.Range("C97").Value = 62
$ws.Range("D97").Value = 9235
$ws.Range("E97").Value = 668

# Row 102
$ws.Range("B102").Value = 8750
$ws.Range("C102").Value = 25
$ws.Range("E102").Value = 911

# Row 110
$ws.Range("B110").Value = 5860
$ws.Range("C110").Value = 92
$ws.Range("D110").Value = 3220
$ws.Range("E110").Value = 2602

# Row 126
$ws.Range("B126").Value = 3954
$ws.Range("C126").Value = 123
$ws.Range("D126").Value = 2844
$ws.Range("E126").Value = 975

# Row 144
$ws.Range("B144").Value = 2756
$ws.Range("C144").Value = 36
$ws.Range("D144").Value = 2318
$ws.Range("E144").Value = 374

# Row 214
$ws.Range("D214").Value = 13
$ws.Range("H214").Value = 0

# Row 215
$ws.Range("D215").Value = 12
$ws.Range("H215").Value = 1
